$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -4
    3  = 0
    7  = 1
    11 = -1
    13 = 1
    17 = 4
    25 = 2
    30 = -5
    38 = 0
    40 = -2
    51 = 4
    55 = -1
    56 = -1
    57 = 3
    58 = 1
    61 = 0
    64 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
